# Fill in the blank "Support Vector" accuracy cell in the results table with
# 66.6412, matching the run-split formatting ("6" hinted east-Asia + rest)
# used throughout the rest of the table.

$d = $word.ActiveDocument

# Locate the single data table and find the one cell whose text is empty
# (just the end-of-cell mark) -- this is robust to row/column reordering.
# Prefer an empty cell that lives in the "Accuracy (%)" row, but fall back
# to any empty cell in the table if that row can't be identified.
$table = $d.Tables.Item(1)
$targetCell = $null
$fallbackCell = $null

for ($r = 1; $r -le $table.Rows.Count; $r++) {
    $rowText = $table.Rows.Item($r).Range.Text
    $isAccuracyRow = $rowText.Contains("Accuracy")
    for ($c = 1; $c -le $table.Columns.Count; $c++) {
        $cell = $table.Cell($r, $c)
        $cellText = $cell.Range.Text
        # An empty cell's Range.Text is just the 2-char end-of-cell mark
        # (CR + cell-mark); any populated cell is longer than that.
        if ($cellText.Length -eq 2) {
            $fallbackCell = $cell
            if ($isAccuracyRow) {
                $targetCell = $cell
            }
        }
    }
}

if ($targetCell -eq $null) {
    $targetCell = $fallbackCell
}

if ($targetCell -eq $null) {
    throw "Could not find the empty accuracy cell to fill in"
}

# Build a minimal single-paragraph OOXML fragment reproducing the same
# run split (first character carrying an eastAsia font hint, remainder in a
# second plain run) seen on every other populated cell in this table, e.g.
# "9" + "8.9724" for the NLP column.
$xml = '<?xml version="1.0" standalone="yes"?>' +
       '<?mso-application progid="Word.Document"?>' +
       '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body>' +
       '<w:p>' +
       '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>6</w:t></w:r>' +
       '<w:r><w:t>6.6412</w:t></w:r>' +
       '</w:p>' +
       '</w:body>' +
       '</w:document>' +
       '</pkg:xmlData>' +
       '</pkg:part>' +
       '</pkg:package>'

$targetRange = $targetCell.Range
$targetRange.Collapse(0)
$null = $targetRange.InsertXML($xml)
